$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1353.0526
$ws.Range("I19").Value = 1221.5714
$ws.Range("J19").Value = 1721.2
$ws.Range("K19").Value = 1221.5714
$ws.Range("L19").Value = 1721.2
$ws.Range("M19").Value = -1046.5714
$ws.Range("N19").Value = -2071.2
$ws.Range("H55").Value = 1562.375
$ws.Range("I55").Value = 2279.8
$ws.Range("J55").Value = 366.66666
$ws.Range("K55").Value = 2279.8
$ws.Range("L55").Value = 366.66666
$ws.Range("M55").Value = -2065.8
$ws.Range("N55").Value = -794.66666
$ws.Range("H76").Value = 4548567.5
$ws.Range("I76").Value = 4765047.5
$ws.Range("J76").Value = 2490
$ws.Range("K76").Value = 4765047.5
$ws.Range("L76").Value = 2490
$ws.Range("M76").Value = -4764732.5
$ws.Range("N76").Value = -3120
$ws.Range("H79").Value = 4548567.5
$ws.Range("I79").Value = 4765047.5
$ws.Range("J79").Value = 2490
$ws.Range("K79").Value = 4765047.5
$ws.Range("L79").Value = 2490
$ws.Range("M79").Value = -4763955.5
$ws.Range("N79").Value = -4674
$ws.Range("H106").Value = 3771.7334
$ws.Range("I106").Value = 3661.4546
$ws.Range("J106").Value = 4075
$ws.Range("K106").Value = 3661.4546
$ws.Range("L106").Value = 4075
$ws.Range("M106").Value = -3030.4546
$ws.Range("N106").Value = -5337
$ws.Range("H132").Value = 2612.2632
$ws.Range("I132").Value = 2758.3125
$ws.Range("J132").Value = 1833.3334
$ws.Range("K132").Value = 8274.9375
$ws.Range("L132").Value = 5500.0002
$ws.Range("M132").Value = -5744.9375
$ws.Range("N132").Value = -10560.0002
$ws.Range("H133").Value = 41433.207
$ws.Range("J133").Value = 41433.207
$ws.Range("L133").Value = 41433.207
$ws.Range("N133").Value = -51553.207
$ws.Range("H134").Value = 44746.5
$ws.Range("J134").Value = 44746.5
$ws.Range("L134").Value = 44746.5
$ws.Range("N134").Value = -54886.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13149.484
$ws.Range("I32").Value = 7852.0527
$ws.Range("J32").Value = 56285.715
$ws.Range("K32").Value = 7852.0527
$ws.Range("L32").Value = 56285.715
$ws.Range("M32").Value = -7565.0527
$ws.Range("N32").Value = -56859.715
$ws.Range("H41").Value = 17333.334
$ws.Range("H110").Value = 901.3333
$ws.Range("I110").Value = 803.7143
$ws.Range("J110").Value = 1038
$ws.Range("K110").Value = 803.7143
$ws.Range("L110").Value = 1038
$ws.Range("M110").Value = 1241.2857
$ws.Range("N110").Value = -5128

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 4227.6665
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()
$ws.Range("H82").Value = 15208
$ws.Range("J82").Value = 19280.666
$ws.Range("L82").Value = 19280.666
$ws.Range("N82").Value = -20046.666
$ws.Range("H85").Value = 15208
$ws.Range("J85").Value = 19280.666
$ws.Range("L85").Value = 19280.666
$ws.Range("N85").Value = -21932.666
$ws.Range("H97").Value = 25545.6
$ws.Range("I97").Value = 964
$ws.Range("J97").Value = 41933.332
$ws.Range("K97").Value = 964
$ws.Range("L97").Value = 41933.332
$ws.Range("M97").Value = 27
$ws.Range("N97").Value = -43915.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2709.8333
$ws.Range("I31").Value = 3076.375
$ws.Range("J31").Value = 1976.75
$ws.Range("K31").Value = 3076.375
$ws.Range("L31").Value = 1976.75
$ws.Range("M31").Value = -2781.375
$ws.Range("N31").Value = -2566.75
$ws.Range("H34").Value = 2709.8333
$ws.Range("I34").Value = 3076.375
$ws.Range("J34").Value = 1976.75
$ws.Range("K34").Value = 3076.375
$ws.Range("L34").Value = 1976.75
$ws.Range("M34").Value = -2874.375
$ws.Range("N34").Value = -2380.75
$ws.Range("H58").Value = 8053.7334
$ws.Range("I58").Value = 1382.6
$ws.Range("J58").Value = 21396
$ws.Range("K58").Value = 1382.6
$ws.Range("L58").Value = 21396
$ws.Range("M58").Value = -1179.6
$ws.Range("N58").Value = -21802
$ws.Range("H99").Value = 3036.1
$ws.Range("I99").Value = 2819.3462
$ws.Range("J99").Value = 4445
$ws.Range("K99").Value = 2819.3462
$ws.Range("L99").Value = 4445
$ws.Range("M99").Value = -1321.3462
$ws.Range("N99").Value = -7441
$ws.Range("H126").Value = 3036.1
$ws.Range("I126").Value = 2819.3462
$ws.Range("J126").Value = 4445
$ws.Range("K126").Value = 8458.0386
$ws.Range("L126").Value = 13335
$ws.Range("M126").Value = -5988.0386
$ws.Range("N126").Value = -18275
$ws.Range("H136").Value = 8053.7334
$ws.Range("I136").Value = 1382.6
$ws.Range("J136").Value = 21396
$ws.Range("K136").Value = 4147.799999999999
$ws.Range("L136").Value = 64188
$ws.Range("M136").Value = -1597.799999999999
$ws.Range("N136").Value = -69288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 512.5
$ws.Range("I32").Value = 450
$ws.Range("J32").Value = 533.3333
$ws.Range("K32").Value = 1350
$ws.Range("L32").Value = 1599.9999
$ws.Range("M32").Value = -1067
$ws.Range("N32").Value = -2165.9999
$ws.Range("H39").Value = 4399.6
$ws.Range("J39").Value = 4399.6
$ws.Range("L39").Value = 13198.8
$ws.Range("N39").Value = -13786.8
$ws.Range("H46").Value = 2344.625
$ws.Range("I46").Value = 2031.4
$ws.Range("J46").Value = 2866.6667
$ws.Range("K46").Value = 6094.200000000001
$ws.Range("L46").Value = 8600.000100000001
$ws.Range("M46").Value = -6003.200000000001
$ws.Range("N46").Value = -8782.000100000001
$ws.Range("H58").Value = 2452.647
$ws.Range("I58").Value = 931.6667
$ws.Range("J58").Value = 2778.5715
$ws.Range("K58").Value = 2795.0001
$ws.Range("L58").Value = 8335.7145
$ws.Range("M58").Value = -2667.0001
$ws.Range("N58").Value = -8591.7145
$ws.Range("H141").Value = 4213.0625
$ws.Range("I141").Value = 1587.1818
$ws.Range("K141").Value = 4761.5454
$ws.Range("M141").Value = 418.4546

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4190.207
$ws.Range("I70").Value = 4096.769
$ws.Range("K70").Value = 4096.769
$ws.Range("M70").Value = -3826.769
$ws.Range("H73").Value = 4190.207
$ws.Range("I73").Value = 4096.769
$ws.Range("K73").Value = 4096.769
$ws.Range("M73").Value = -3160.769
$ws.Range("H92").Value = 9387.5
$ws.Range("J92").Value = 9387.5
$ws.Range("L92").Value = 9387.5
$ws.Range("N92").Value = -13131.5
$ws.Range("H126").Value = 2985.8572
$ws.Range("I126").Value = 3221.1538
$ws.Range("K126").Value = 9663.4614
$ws.Range("M126").Value = -7193.4614
$ws.Range("H132").Value = 2796
$ws.Range("I132").Value = 2437.8215
$ws.Range("J132").Value = 4049.625
$ws.Range("K132").Value = 7313.4645
$ws.Range("L132").Value = 12148.875
$ws.Range("M132").Value = -4783.4645
$ws.Range("N132").Value = -17208.875
$ws.Range("H136").Value = 8834.666999999999
$ws.Range("J136").Value = 8834.666999999999
$ws.Range("L136").Value = 26504.001
$ws.Range("N136").Value = -31604.001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1012241.9
$ws.Range("I40").Value = 1264739.2
$ws.Range("J40").Value = 2252.5
$ws.Range("K40").Value = 1264739.2
$ws.Range("L40").Value = 2252.5
$ws.Range("M40").Value = -1264603.2
$ws.Range("N40").Value = -2524.5
$ws.Range("H136").Value = 4576.6
$ws.Range("I136").Value = 1672.6666
$ws.Range("J136").Value = 8932.5
$ws.Range("K136").Value = 5017.9998
$ws.Range("L136").Value = 26797.5
$ws.Range("M136").Value = -2467.9998
$ws.Range("N136").Value = -31897.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 16690
$ws.Range("I51").Value = 5035
$ws.Range("K51").Value = 5035
$ws.Range("M51").Value = -4525
$ws.Range("H126").Value = 770.2222
$ws.Range("I126").Value = 678.375
$ws.Range("J126").Value = 1505
$ws.Range("K126").Value = 2035.125
$ws.Range("L126").Value = 4515
$ws.Range("M126").Value = 434.875
$ws.Range("N126").Value = -9455
$ws.Range("H132").Value = 1772.3182
$ws.Range("I132").Value = 1513.7858
$ws.Range("J132").Value = 2224.75
$ws.Range("K132").Value = 4541.357400000001
$ws.Range("L132").Value = 6674.25
$ws.Range("M132").Value = -2011.357400000001
$ws.Range("N132").Value = -11734.25
$ws.Range("H136").Value = 706.7727
$ws.Range("I136").Value = 665.3125
$ws.Range("J136").Value = 817.3333
$ws.Range("K136").Value = 1995.9375
$ws.Range("L136").Value = 2451.9999
$ws.Range("M136").Value = 554.0625
$ws.Range("N136").Value = -7551.9999
